$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1367.1765
$ws.Range("I19").Value = 726.2
$ws.Range("K19").Value = 726.2
$ws.Range("M19").Value = -551.2
$ws.Range("H98").Value = 877.36
$ws.Range("I98").Value = 830.375
$ws.Range("K98").Value = 830.375
$ws.Range("M98").Value = 667.625
$ws.Range("H122").Value = 877.36
$ws.Range("I122").Value = 830.375
$ws.Range("K122").Value = 2491.125
$ws.Range("M122").Value = -41.125
$ws.Range("H123").Value = 38000
$ws.Range("J123").Value = 38000
$ws.Range("L123").Value = 38000
$ws.Range("N123").Value = -47800
$ws.Range("H138").Value = 2177.79
$ws.Range("J138").Value = 2328.191
$ws.Range("L138").Value = 6984.572999999999
$ws.Range("N138").Value = -17264.573

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4300.2163
$ws.Range("I32").Value = 4591.5
$ws.Range("K32").Value = 4591.5
$ws.Range("M32").Value = -4304.5
$ws.Range("H45").Value = 2438.647
$ws.Range("I45").Value = 2195.7
$ws.Range("J45").Value = 2785.7144
$ws.Range("K45").Value = 2195.7
$ws.Range("L45").Value = 2785.7144
$ws.Range("M45").Value = -1818.7
$ws.Range("N45").Value = -3539.7144
$ws.Range("H61").Value = 966.3158
$ws.Range("I61").Value = 966.3158
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 966.3158
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -754.3158
$ws.Range("N61").Value = -754.3158
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").Value = 0
$ws.Range("H132").Value = 3009.6206
$ws.Range("I132").Value = 2845.0833
$ws.Range("K132").Value = 8535.249899999999
$ws.Range("M132").Value = -6005.249899999999
$ws.Range("H136").Value = 966.3158
$ws.Range("I136").Value = 966.3158
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 2898.9474
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -348.9474
$ws.Range("N136").Value = -348.9474

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4048.0588
$ws.Range("I86").Value = 3876.0625
$ws.Range("K86").Value = 3876.0625
$ws.Range("M86").Value = -2753.0625
$ws.Range("H89").Value = 4048.0588
$ws.Range("I89").Value = 3876.0625
$ws.Range("K89").Value = 19380.3125
$ws.Range("M89").Value = -13764.3125
$ws.Range("H99").Value = 55556696
$ws.Range("I99").Value = 71429590
$ws.Range("J99").Value = 1550
$ws.Range("K99").Value = 71429590
$ws.Range("L99").Value = 1550
$ws.Range("M99").Value = -71428092
$ws.Range("N99").Value = -4546
$ws.Range("H134").Value = 9641.25
$ws.Range("I134").Value = 5971
$ws.Range("K134").Value = 17913
$ws.Range("M134").Value = -15378

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 724.561
$ws.Range("I31").Value = 651.125
$ws.Range("J31").Value = 882.7308
$ws.Range("K31").Value = 651.125
$ws.Range("L31").Value = 882.7308
$ws.Range("M31").Value = -356.125
$ws.Range("N31").Value = -1472.7308
$ws.Range("H34").Value = 724.561
$ws.Range("I34").Value = 651.125
$ws.Range("J34").Value = 882.7308
$ws.Range("K34").Value = 651.125
$ws.Range("L34").Value = 882.7308
$ws.Range("M34").Value = -449.125
$ws.Range("N34").Value = -1286.7308
$ws.Range("H63").Value = 30000
$ws.Range("J63").Value = 30000
$ws.Range("L63").Value = 30000
$ws.Range("N63").Value = -31372
$ws.Range("H66").Value = 30000
$ws.Range("J66").Value = 30000
$ws.Range("L66").Value = 90000
$ws.Range("N66").Value = -96864
$ws.Range("H132").Value = 7317.0557
$ws.Range("I132").Value = 7980.533
$ws.Range("J132").Value = 3999.6667
$ws.Range("K132").Value = 23941.599
$ws.Range("L132").Value = 11999.0001
$ws.Range("M132").Value = -21411.599
$ws.Range("N132").Value = -17059.0001
$ws.Range("H134").Value = 8548141
$ws.Range("I134").Value = 9524760
$ws.Range("K134").Value = 28574280
$ws.Range("M134").Value = -28571745

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1336.5294
$ws.Range("I5").Value = 1209.52
$ws.Range("J5").Value = 1689.3334
$ws.Range("K5").Value = 3628.56
$ws.Range("L5").Value = 5068.0002
$ws.Range("M5").Value = -3516.56
$ws.Range("N5").Value = -5292.0002
$ws.Range("H131").Value = 33335016
$ws.Range("J131").Value = 1977.1305
$ws.Range("L131").Value = 5931.3915
$ws.Range("N131").Value = -16011.3915
$ws.Range("H135").Value = 1336.5294
$ws.Range("I135").Value = 1209.52
$ws.Range("J135").Value = 1689.3334
$ws.Range("K135").Value = 10885.68
$ws.Range("L135").Value = 15204.0006
$ws.Range("M135").Value = -8350.68
$ws.Range("N135").Value = -20274.0006

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 10015476
$ws.Range("I11").Value = 10643421
$ws.Range("J11").Value = 4050000
$ws.Range("K11").Value = 10643421
$ws.Range("L11").Value = 4050000
$ws.Range("M11").Value = -10643282
$ws.Range("N11").Value = -4050278
$ws.Range("H45").Value = 37181.25
$ws.Range("J45").Value = 37181.25
$ws.Range("L45").Value = 37181.25
$ws.Range("N45").Value = -38299.25
$ws.Range("H70").Value = 14520396
$ws.Range("I70").Value = 20837184
$ws.Range("J70").Value = 10530846
$ws.Range("K70").Value = 20837184
$ws.Range("L70").Value = 10530846
$ws.Range("M70").Value = -20836914
$ws.Range("N70").Value = -10531386
$ws.Range("H73").Value = 14520396
$ws.Range("I73").Value = 20837184
$ws.Range("J73").Value = 10530846
$ws.Range("K73").Value = 20837184
$ws.Range("L73").Value = 10530846
$ws.Range("M73").Value = -20836248
$ws.Range("N73").Value = -10532718
$ws.Range("H124").Value = 46966.668
$ws.Range("J124").Value = 46966.668
$ws.Range("L124").Value = 46966.668
$ws.Range("N124").Value = -56786.668
$ws.Range("H132").Value = 1894.1428
$ws.Range("I132").Value = 1544.9
$ws.Range("J132").Value = 3989.6
$ws.Range("K132").Value = 4634.700000000001
$ws.Range("L132").Value = 11968.8
$ws.Range("M132").Value = -2104.700000000001
$ws.Range("N132").Value = -17028.8

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2017.1177
$ws.Range("I7").Value = 1839.6
$ws.Range("J7").Value = 2270.7144
$ws.Range("K7").Value = 1839.6
$ws.Range("L7").Value = 2270.7144
$ws.Range("M7").Value = -1727.6
$ws.Range("N7").Value = -2494.7144
$ws.Range("H25").Value = 713336
$ws.Range("J25").Value = 70004
$ws.Range("L25").Value = 70004
$ws.Range("N25").Value = -70464
$ws.Range("H40").Value = 3630
$ws.Range("I40").Value = 2890
$ws.Range("J40").Value = 4000
$ws.Range("K40").Value = 2890
$ws.Range("L40").Value = 4000
$ws.Range("M40").Value = -2754
$ws.Range("N40").Value = -4272
$ws.Range("H55").Value = 499.17648
$ws.Range("I55").Value = 50.57143
$ws.Range("J55").Value = 813.2
$ws.Range("K55").Value = 50.57143
$ws.Range("L55").Value = 813.2
$ws.Range("M55").Value = 122.42857
$ws.Range("N55").Value = -1159.2
$ws.Range("H62").Value = 14500
$ws.Range("J62").Value = 14500
$ws.Range("L62").Value = 14500
$ws.Range("N62").Value = -15748
$ws.Range("H65").Value = 14500
$ws.Range("J65").Value = 14500
$ws.Range("L65").Value = 43500
$ws.Range("N65").Value = -49740
$ws.Range("H68").Value = 1354.2
$ws.Range("I68").Value = 1193.8182
$ws.Range("J68").Value = 1795.25
$ws.Range("K68").Value = 1193.8182
$ws.Range("L68").Value = 1795.25
$ws.Range("M68").Value = -444.8181999999999
$ws.Range("N68").Value = -3293.25
$ws.Range("H71").Value = 1354.2
$ws.Range("I71").Value = 1193.8182
$ws.Range("J71").Value = 1795.25
$ws.Range("K71").Value = 5969.090999999999
$ws.Range("L71").Value = 8976.25
$ws.Range("M71").Value = -2225.090999999999
$ws.Range("N71").Value = -16464.25
$ws.Range("H82").Value = 1840.7142
$ws.Range("I82").Value = 1697.5
$ws.Range("K82").Value = 1697.5
$ws.Range("M82").Value = -1336.5
$ws.Range("H85").Value = 1840.7142
$ws.Range("I85").Value = 1697.5
$ws.Range("K85").Value = 1697.5
$ws.Range("M85").Value = -449.5
$ws.Range("H122").Value = 31252264
$ws.Range("J122").Value = 2451.25
$ws.Range("L122").Value = 7353.75
$ws.Range("N122").Value = -12253.75
$ws.Range("H126").Value = 2017.1177
$ws.Range("I126").Value = 1839.6
$ws.Range("J126").Value = 2270.7144
$ws.Range("K126").Value = 5518.799999999999
$ws.Range("L126").Value = 6812.1432
$ws.Range("M126").Value = -3048.799999999999
$ws.Range("N126").Value = -11752.1432
$ws.Range("H132").Value = 65337.625
$ws.Range("I132").Value = 2933.1667
$ws.Range("J132").Value = 102780.3
$ws.Range("K132").Value = 8799.500100000001
$ws.Range("L132").Value = 308340.9
$ws.Range("M132").Value = -6269.500100000001
$ws.Range("N132").Value = -313400.9
$ws.Range("H136").Value = 1925
$ws.Range("I136").Value = 1400
$ws.Range("K136").Value = 4200
$ws.Range("M136").Value = -1650

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 117413.664
$ws.Range("J14").Value = 895.4
$ws.Range("L14").Value = 895.4
$ws.Range("N14").Value = -1231.4
$ws.Range("H122").Value = 18573656
$ws.Range("I122").Value = 28891344
$ws.Range("J122").Value = 1820
$ws.Range("K122").Value = 86674032
$ws.Range("L122").Value = 5460
$ws.Range("M122").Value = -86671582
$ws.Range("N122").Value = -10360
$ws.Range("H132").Value = 1872.3103
$ws.Range("I132").Value = 1595.3636
$ws.Range("J132").Value = 2742.7144
$ws.Range("K132").Value = 4786.0908
$ws.Range("L132").Value = 8228.143199999999
$ws.Range("M132").Value = -2256.0908
$ws.Range("N132").Value = -13288.1432
$ws.Range("H136").Value = 1624
$ws.Range("I136").Value = 1200.4286
$ws.Range("J136").Value = 2118.1667
$ws.Range("K136").Value = 3601.2858
$ws.Range("L136").Value = 6354.500100000001
$ws.Range("M136").Value = -1051.2858
$ws.Range("N136").Value = -11454.5001
